# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Update DAMSLTag (col I) and DialogAct (col J)
# values for the rows whose re-annotation changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 10;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 12;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 17;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 24;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 34;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 45;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 52;  I = "%";  J = "Uninterpretable" },
    @{ Row = 62;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 79;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 81;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 82;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 108; I = "aa"; J = "Agree/Accept" },
    @{ Row = 110; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 113; I = "%";  J = "Uninterpretable" },
    @{ Row = 123; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 126; I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 131; I = "%";  J = "Uninterpretable" },
    @{ Row = 139; I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 140; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 150; I = "aa"; J = "Agree/Accept" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
